$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (most recent period) balance-sheet figures for JNJ
$ws.Range("B2").Value = 24619000000.0
$ws.Range("B3").Value = 14938000000.0
$ws.Range("B4").Value = 9952000000.0
$ws.Range("B6").Value = 52533000000.0
$ws.Range("B7").Value = 18367000000.0
$ws.Range("B9").Value = 86798000000.0
$ws.Range("B10").Value = 8321000000.0
$ws.Range("B11").Value = 6538000000.0
$ws.Range("B12").Value = 120024000000.0
$ws.Range("B13").Value = 172557000000.0
$ws.Range("B14").Value = 3350000000.0
$ws.Range("B15").Value = 8503000000.0
$ws.Range("B16").Value = 27202000000.0
$ws.Range("B18").Value = 1877000000.0
$ws.Range("B20").Value = 40932000000.0
$ws.Range("B21").Value = 30263000000.0
$ws.Range("B22").Value = 10512000000.0
$ws.Range("B23").Value = 6507000000.0
$ws.Range("B24").Value = 18509000000.0
$ws.Range("B25").Value = 65791000000.0
$ws.Range("B26").Value = 106723000000.0
$ws.Range("B27").Value = 3120000000.0
$ws.Range("B28").Value = 116508000000.0
$ws.Range("B29").Value = 38466000000.0
$ws.Range("B30").Value = 65834000000.0
$ws.Range("B31").Value = 65834000000.0
$ws.Range("B32").Value = 172557000000.0
$ws.Range("B33").Value = 2632702000.0
$ws.Range("B34").Value = -20964000000.0

# Widen column B to match the other data columns
$ws.Columns("B").ColumnWidth = 16.885714285714286
